# Update threshold values on Sheet1:
#   - beta_distance_range  (row 3): Min 4.5 -> 5,  Max 9.3 -> 9
#   - alpha_distance_range (row 2): Max 10.6 -> 10
#   - ratio_threshold_range(row 4): Max 1.5 -> 1.4
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 10
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 9
$ws.Range("C4").Value = 1.4
